$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @(45923.01041666666,45923.02083333334,45923.03125,45923.04166666666,45923.05208333334,45923.0625,45923.07291666666,45923.08333333334,45923.09375,45923.10416666666,45923.11458333334,45923.125,45923.13541666666,45923.14583333334,45923.15625,45923.16666666666,45923.17708333334,45923.1875,45923.19791666666,45923.20833333334,45923.21875,45923.22916666666,45923.23958333334,45923.25,45923.26041666666,45923.27083333334,45923.28125,45923.29166666666,45923.30208333334,45923.3125,45923.32291666666,45923.33333333334,45923.34375,45923.35416666666,45923.36458333334,45923.375,45923.38541666666,45923.39583333334,45923.40625,45923.41666666666,45923.42708333334,45923.4375,45923.44791666666,45923.45833333334,45923.46875,45923.47916666666,45923.48958333334,45923.5,45923.51041666666,45923.52083333334,45923.53125,45923.54166666666,45923.55208333334,45923.5625,45923.57291666666,45923.58333333334,45923.59375,45923.60416666666,45923.61458333334,45923.625,45923.63541666666,45923.64583333334,45923.65625,45923.66666666666,45923.67708333334,45923.6875,45923.69791666666,45923.70833333334,45923.71875,45923.72916666666,45923.73958333334,45923.75,45923.76041666666,45923.77083333334,45923.78125,45923.79166666666,45923.80208333334,45923.8125,45923.82291666666,45923.83333333334,45923.84375,45923.85416666666,45923.86458333334,45923.875,45923.88541666666,45923.89583333334,45923.90625,45923.91666666666,45923.92708333334,45923.9375,45923.94791666666,45923.95833333334,45923.96875,45923.97916666666,45923.98958333334,45924)
$values = @(605,602,599,595,467,464,461,459,382,380,379,378,307,306,306,305,265,264,264,264,236,236,236,235,225,225,225,225,212,212,212,212,195,197,198,200,205,206,207,209,207,207,207,207,197,197,197,196,194,194,193,193,193,193,193,193,208,208,208,208,227,227,226,225,279,279,279,279,361,363,364,366,447,448,449,450,461,460,459,459,408,407,407,407,347,346,346,346,295,295,294,294,0,0,0,0)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $timestamps[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
